$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row of quotations (2025-09-13) below the existing data.
$newRow = 9

# Column A: date value, formatted like the cell above it (A8).
$ws.Cells.Item($newRow, 1).Value = 45913
$ws.Cells.Item($newRow, 1).NumberFormat = $ws.Cells.Item($newRow - 1, 1).NumberFormat

# Columns B..E: textual quotation values (stored as text, Portuguese decimal comma).
$ws.Cells.Item($newRow, 2).Value = "21,1936"
$ws.Cells.Item($newRow, 3).Value = "14,9727"
$ws.Cells.Item($newRow, 4).Value = "14,8984"
$ws.Cells.Item($newRow, 5).Value = "14,8984"
